# JekyllAndHyde.xlsx re-upload:
#   - "jekyll" sheet becomes the active/selected sheet (was "Formatted")
#   - On "jekyll", columns A ("Shuffle") and B ("sentenceID") are unhidden
#     and given explicit widths instead of being hidden with width 0
#   - Column A still holds the volatile =RAND() shuffle key, so its cached
#     values simply reflect whatever the engine recalculates on save

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("jekyll")

# Unhide column A ("Shuffle") and give it its authored width
$colA = $ws1.Columns.Item(1)
$colA.Hidden = $false
$colA.ColumnWidth = 11

# Unhide column B ("sentenceID") and give it its authored width
$colB = $ws1.Columns.Item(2)
$colB.Hidden = $false
$colB.ColumnWidth = 11.85

# Make "jekyll" the active/selected sheet (previously "Formatted" was
# the selected tab) - this flips tabSelected on the sheetViews and
# updates the workbook's bookView activeTab accordingly.
$ws1.Activate()
